$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Id value in K2 from 0247 to 0263.
# The leading apostrophe forces Excel to keep the numeric-looking value
# stored as text (matching the cell's existing "quote prefix" text style),
# instead of converting it to a real number and dropping the leading zero.
$ws.Range("K2").Value = "'0263"

# Update the active selection to K2 (matches the last edited cell)
$ws.Range("K2").Select()
